$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.036.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.385.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.383.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("E11").Value = "  -2.44%  "
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.963.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("E16").Value = "  -2.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.381.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.130.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("E19").Value = "  -1.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.526.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  +10.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.174"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.90%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -2.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.774"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.471.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("E51").Value = "  +4.06%  "
